$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.219.27'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').Value = '2.216.61'
$ws.Range('E3').Value = '  -1.02%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.41'
$ws.Range('E5').Value = '  -1.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '88.58'
$ws.Range('E6').Value = '  -5.23%  '
$ws.Range('E7').Value = '  -3.21%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.491'
$ws.Range('E9').Value = '  -5.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.39'
$ws.Range('E10').Value = '  -2.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0778'
$ws.Range('E11').Value = '  -3.52%  '
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.86'
$ws.Range('E13').Value = '  -3.76%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.554.51'
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.289.56'
$ws.Range('E15').Value = '  -2.06%  '
$ws.Range('E16').Value = '  -3.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.08'
$ws.Range('E17').Value = '  -2.39%  '
$ws.Range('D18').Value = '43.938.07'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = '0.0₃0900'
$ws.Range('E19').Value = '  -6.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.95'
$ws.Range('E20').Value = '  -5.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.16'
$ws.Range('E21').Value = '  -7.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '64.07'
$ws.Range('E22').Value = '  -2.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.19'
$ws.Range('E23').Value = '  -1.91%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.85'
$ws.Range('E24').Value = '  -8.05%  '
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('E26').Value = '  -5.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.24'
$ws.Range('E27').Value = '  +1.83%  '
$ws.Range('E28').Value = '  -4.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.04'
$ws.Range('E29').Value = '  -8.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.35'
$ws.Range('E30').Value = '  -3.35%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.57'
$ws.Range('E31').Value = '  -5.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '145.89'
$ws.Range('E32').Value = '  -3.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.59'
$ws.Range('E33').Value = '  -0.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0752'
$ws.Range('E34').Value = '  -5.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.92'
$ws.Range('E35').Value = '  -4.70%  '
$ws.Range('E36').Value = '  -2.89%  '
$ws.Range('E37').Value = '  -4.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.70'
$ws.Range('E38').Value = '  -2.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.13'
$ws.Range('E39').Value = '  -0.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.18'
$ws.Range('E40').Value = '  -8.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.60'
$ws.Range('E41').Value = '  -4.51%  '
$ws.Range('E42').Value = '  -4.11%  '
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').Value = '1.735.33'
$ws.Range('E44').Value = '  +1.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.66'
$ws.Range('E45').Value = '  +4.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '77.99'
$ws.Range('E46').Value = '  -5.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.178'
$ws.Range('E47').Value = '  -6.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '94.35'
$ws.Range('E48').Value = '  -4.91%  '
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '66.25'
$ws.Range('E49').Value = '  -1.86%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.60'
$ws.Range('E50').Value = '  -6.32%  '
$ws.Range('D51').Value = '2.435.40'
$ws.Range('E51').Value = '  -1.00%  '
